$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$newData = @(
    @(465, 44539, 5, 31, 194.1139636819036),
    @(466, 44540, 4, 23, 144.0200375704446),
    @(467, 44541, 0, 22, 137.7582968065122),
    @(468, 44542, 14, 36, 225.4226675015654),
    @(469, 44543, 7, 33, 206.6374452097683),
    @(470, 44544, 8, 41, 256.7313713212273),
    @(471, 44545, 0, 38, 237.9461490294302),
    @(472, 44546, 7, 40, 250.4696305572949),
    @(473, 44547, 16, 52, 325.6105197244834),
    @(474, 44548, 4, 56, 350.6574827802129),
    @(475, 44550, 17, 59, 369.44270507201),
    @(476, 44551, 22, 74, 463.3688165309956),
    @(477, 44552, 1, 67, 419.536631183469),
    @(478, 44553, 4, 71, 444.5835942391985),
    @(479, 44554, 9, 73, 457.1070757670633),
    @(480, 44555, 20, 77, 482.1540388227927),
    @(481, 44556, 2, 75, 469.6305572949279),
    @(482, 44557, 24, 82, 513.4627426424546),
    @(483, 44558, 9, 69, 432.0601127113338),
    @(484, 44559, 8, 76, 475.8922980588604),
    @(485, 44560, 11, 83, 519.724483406387),
    @(486, 44561, 31, 105, 657.4827802128992),
    @(487, 44562, 29, 114, 713.8384470882905),
    @(488, 44563, 6, 118, 738.8854101440201),
    @(489, 44564, 27, 121, 757.6706324358171),
    @(490, 44565, 28, 140, 876.6437069505322),
    @(491, 44566, 25, 157, 983.0932999373825)
)

$lastRow = 464
$firstNewRow = 465
$lastNewRow = 491

# Extend formatting (column A date style, border, alignment) down from the last
# existing row to cover all the freshly appended rows.
$ws.Range("A$lastRow`:D$lastRow").Copy()
$ws.Range("A$firstNewRow`:D$lastNewRow").PasteSpecial(-4122)

foreach ($row in $newData) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value2 = $row[1]
    $ws.Cells.Item($r, 2).Value2 = $row[2]
    $ws.Cells.Item($r, 3).Value2 = $row[3]
    $ws.Cells.Item($r, 4).Value2 = $row[4]
}
